$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete cell C2 entirely (remove value, making the cell empty so the <c> element disappears)
$ws.Range("C2").ClearContents()

# Update E2
$ws.Range("E2").Value = 8.045645122021927

# Update C3..C19, E13, E15, E16, E19 per diff
$ws.Range("C3").Value = -5.232639093663827
$ws.Range("C4").Value = 2.124540184802992
$ws.Range("C5").Value = 8.081020954067753
$ws.Range("C6").Value = 4.489210662380971
$ws.Range("C7").Value = -0.8752093743685241
$ws.Range("C9").Value = 4.818339085077561
$ws.Range("C11").Value = 4.613634856640769
$ws.Range("E13").Value = 4.838485897465628
$ws.Range("E15").Value = 0.869978169785246
$ws.Range("E16").Value = 2.551560717335266
$ws.Range("C18").Value = -2.447533648174649
$ws.Range("C19").Value = 1.038949519463617
$ws.Range("E19").Value = -1.220869074712128
